# Update crypto price/volume figures per the Fri Oct 11 23:33:49 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.511.56"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "2.406.65"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'573.11"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "'145.75"
$ws.Range("E6").Value = "  +4.95%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "2.432.86"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +4.60%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "'27.43"
$ws.Range("E14").Value = "  +6.09%  "
$ws.Range("D16").Value = "2.883.71"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "61.792.43"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "2.422.88"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "'7.87"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").Value = "'10.93"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("D21").Value = "'327.60"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  +12.00%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'65.53"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("B26").Value = "BabyDogeCoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D26").Value = "0.0₆0585"
$ws.Range("E26").Value = "  +100.52%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'624.68"
$ws.Range("E27").Value = "  +10.54%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.10"
$ws.Range("E28").Value = "  +10.54%  "
$ws.Range("D29").Value = "'8.55"
$ws.Range("E29").Value = "  +5.30%  "
$ws.Range("D30").Value = "0.0₃0981"
$ws.Range("E30").Value = "  +5.48%  "
$ws.Range("D31").Value = "2.561.12"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("E33").Value = "  +7.02%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.137"
$ws.Range("E34").Value = "  +3.63%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "'1.85"
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "'151.90"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("D42").Value = "'18.60"
$ws.Range("E43").Value = "  +13.37%  "
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'14.75"
$ws.Range("D47").Value = "'144.27"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").Value = "'20.51"
$ws.Range("E49").Value = "  +6.74%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("E51").Value = "  +2.46%  "
